$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Förändrad) rows 2-16 from 2023-10-22 (45221) to 2023-10-25 (45224)
for ($r = 2; $r -le 16; $r++) {
    $ws.Cells.Item($r, 3).Value = 45224
}
